$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update D3 value (shared string index 11 -> 14, i.e. "10*10^4" -> "9*10^4")
$ws.Range("D3").Value = "9*10^4"

# Update the selection on the sheet
$ws.Range("E18").Select()

# Update the workbook window position/size
$excel.ActiveWindow.Left = 4560
$excel.ActiveWindow.Top = 2385
$excel.ActiveWindow.Width = 10350
$excel.ActiveWindow.Height = 11055
